# "se controlan puntos de TPO" - mark Marcelo's column (E) as meeting the
# criteria for rows 5, 6, 8 and 9 (same "X" used elsewhere in column G).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E5").Value = "X"
$ws.Range("E6").Value = "X"
$ws.Range("E8").Value = "X"
$ws.Range("E9").Value = "X"

# Reflect the author's new cursor position/selection on the sheet.
$ws.Range("E9").Select()
